# Update workbook to reflect data as of 2022-04-10
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet/tab to reflect the new "through" date
$ws.Name = "Through 2022-04-10"

# Update the label in A5 (April row header) to reflect the new "through" date
$ws.Range("A5").Value = "April (through 04-10)"

# Update April row (row 5) values
$ws.Range("B5").Value = 8
$ws.Range("C5").Value = 9
$ws.Range("D5").Value = 18
$ws.Range("E5").Value = 17
$ws.Range("F5").Value = 16
$ws.Range("G5").Value = 22
$ws.Range("H5").Value = 29
$ws.Range("I5").Value = 38

# Update Total row (row 6) values
$ws.Range("B6").Value = 74
$ws.Range("C6").Value = 137
$ws.Range("D6").Value = 207
$ws.Range("E6").Value = 214
$ws.Range("F6").Value = 126
$ws.Range("G6").Value = 220
$ws.Range("H6").Value = 452
$ws.Range("I6").Value = 472
